# "added 2015 NRI data"
# Appends a new 2015 row (row 9) to the two per-year detail sheets
# (county-year, point-year) and refreshes the aggregated totals on the
# two "-overall" summary sheets (county-overall, point-overall) so they
# include the new year's observations.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# county-year: new row 9 (year 2015)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("county-year")

# Clone the formatting (number format + border) of the last data row so
# the new row matches the existing styled rows.
$ws.Range("A8:AD8").Copy()
$ws.Range("A9:AD9").PasteSpecial(-4122)

$countyYearRow9 = @(2015,0,0,3072,100,0,0,100,0,0,0,0,3072,100,0,0,3072,100,0,0,3072,100,3065,99.772132873535156,7,0.2278645783662796,2828,92.057289123535156,244,7.9427084922790527)
for ($i = 0; $i -lt $countyYearRow9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $countyYearRow9[$i]
}

# ---------------------------------------------------------------------
# point-year: new row 9 (year 2015)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("point-year")

$ws.Range("A8:AD8").Copy()
$ws.Range("A9:AD9").PasteSpecial(-4122)

$pointYearRow9 = @(2015,0,0,1362620,100,0,0,100,0,0,0,0,1362620,100,0,0,1362620,100,0,0,1362620,100,1361164,99.893150329589844,1456,0.10685297101736069,1226047,89.9771728515625,136573,10.022823333740234)
for ($i = 0; $i -lt $pointYearRow9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $pointYearRow9[$i]
}

# ---------------------------------------------------------------------
# county-overall: refresh row 2 totals to include the 2015 observations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("county-overall")

$ws.Range("B2").Value = 25.03662109375
$ws.Range("C2").Value = 24576
$ws.Range("K2").Value = 72.074378967285156
$ws.Range("L2").Value = 6863
$ws.Range("M2").Value = 27.925619125366211
$ws.Range("O2").Value = 31.298828125
$ws.Range("P2").Value = 16884
$ws.Range("Q2").Value = 68.701171875
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = 12288
$ws.Range("U2").Value = 50
$ws.Range("V2").Value = 21455
$ws.Range("W2").Value = 87.300621032714844
$ws.Range("X2").Value = 3121
$ws.Range("Y2").Value = 12.699381828308105
$ws.Range("Z2").Value = 14119
$ws.Range("AA2").Value = 57.450359344482422
$ws.Range("AB2").Value = 10457
$ws.Range("AC2").Value = 42.549640655517578

# ---------------------------------------------------------------------
# point-overall: refresh row 2 totals to include the 2015 observations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("point-overall")

$ws.Range("B2").Value = 25.396230697631836
$ws.Range("C2").Value = 10900960
$ws.Range("K2").Value = 74.413925170898438
$ws.Range("L2").Value = 2789128
$ws.Range("M2").Value = 25.586076736450195
$ws.Range("O2").Value = 30.812211990356445
$ws.Range("P2").Value = 7542133
$ws.Range("Q2").Value = 69.187789916992188
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = 5450480
$ws.Range("U2").Value = 50
$ws.Range("V2").Value = 9528148
$ws.Range("W2").Value = 87.406501770019531
$ws.Range("X2").Value = 1372812
$ws.Range("Y2").Value = 12.593496322631836
$ws.Range("Z2").Value = 6172671
$ws.Range("AA2").Value = 56.625022888183594
$ws.Range("AB2").Value = 4728289
$ws.Range("AC2").Value = 43.374977111816406
